$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header B1: TASA_DESOCUPACION -> TASA_OCUPACION
$ws.Range("B1").Value = "TASA_OCUPACION"

# Replace column B data (rows 2-115) with updated TASA_OCUPACION series
$ws.Cells.Item(2, 2).Value = 59.165300000000002
$ws.Cells.Item(3, 2).Value = 59.96
$ws.Cells.Item(4, 2).Value = 60.731699999999996
$ws.Cells.Item(5, 2).Value = 61.433399999999999
$ws.Cells.Item(6, 2).Value = 61.3508
$ws.Cells.Item(7, 2).Value = 61.589100000000002
$ws.Cells.Item(8, 2).Value = 60.8262
$ws.Cells.Item(9, 2).Value = 61.186100000000003
$ws.Cells.Item(10, 2).Value = 61.272599999999997
$ws.Cells.Item(11, 2).Value = 63.287500000000001
$ws.Cells.Item(12, 2).Value = 62.933799999999998
$ws.Cells.Item(13, 2).Value = 61.516100000000002
$ws.Cells.Item(14, 2).Value = 58.956000000000003
$ws.Cells.Item(15, 2).Value = 59.887900000000002
$ws.Cells.Item(16, 2).Value = 59.194699999999997
$ws.Cells.Item(17, 2).Value = 61.188400000000001
$ws.Cells.Item(18, 2).Value = 60.2331
$ws.Cells.Item(19, 2).Value = 61.089599999999997
$ws.Cells.Item(20, 2).Value = 59.348300000000002
$ws.Cells.Item(21, 2).Value = 60.538600000000002
$ws.Cells.Item(22, 2).Value = 60.475099999999998
$ws.Cells.Item(23, 2).Value = 62.264400000000002
$ws.Cells.Item(24, 2).Value = 61.978999999999999
$ws.Cells.Item(25, 2).Value = 60.709800000000001
$ws.Cells.Item(26, 2).Value = 58.301499999999997
$ws.Cells.Item(27, 2).Value = 59.127099999999999
$ws.Cells.Item(28, 2).Value = 59.236600000000003
$ws.Cells.Item(29, 2).Value = 60.873199999999997
$ws.Cells.Item(30, 2).Value = 59.843800000000002
$ws.Cells.Item(31, 2).Value = 61.321100000000001
$ws.Cells.Item(32, 2).Value = 59.432499999999997
$ws.Cells.Item(33, 2).Value = 60.1233
$ws.Cells.Item(34, 2).Value = 59.701700000000002
$ws.Cells.Item(35, 2).Value = 61.355499999999999
$ws.Cells.Item(36, 2).Value = 60.614800000000002
$ws.Cells.Item(37, 2).Value = 60.211300000000001
$ws.Cells.Item(38, 2).Value = 57.386400000000002
$ws.Cells.Item(39, 2).Value = 57.833300000000001
$ws.Cells.Item(40, 2).Value = 58.7333
$ws.Cells.Item(41, 2).Value = 59.837899999999998
$ws.Cells.Item(42, 2).Value = 59.254100000000001
$ws.Cells.Item(43, 2).Value = 59.891800000000003
$ws.Cells.Item(44, 2).Value = 58.879300000000001
$ws.Cells.Item(45, 2).Value = 59.593800000000002
$ws.Cells.Item(46, 2).Value = 59.470599999999997
$ws.Cells.Item(47, 2).Value = 60.070599999999999
$ws.Cells.Item(48, 2).Value = 58.795200000000001
$ws.Cells.Item(49, 2).Value = 59.779000000000003
$ws.Cells.Item(50, 2).Value = 56.465400000000002
$ws.Cells.Item(51, 2).Value = 57.551000000000002
$ws.Cells.Item(52, 2).Value = 57.6751
$ws.Cells.Item(53, 2).Value = 57.314100000000003
$ws.Cells.Item(54, 2).Value = 57.578800000000001
$ws.Cells.Item(55, 2).Value = 58.783099999999997
$ws.Cells.Item(56, 2).Value = 57.084800000000001
$ws.Cells.Item(57, 2).Value = 56.8598
$ws.Cells.Item(58, 2).Value = 57.3598
$ws.Cells.Item(59, 2).Value = 58.864100000000001
$ws.Cells.Item(60, 2).Value = 58.744900000000001
$ws.Cells.Item(61, 2).Value = 58.192399999999999
$ws.Cells.Item(62, 2).Value = 55.3277
$ws.Cells.Item(63, 2).Value = 56.0715
$ws.Cells.Item(64, 2).Value = 52.45
$ws.Cells.Item(65, 2).Value = 42.497100000000003
$ws.Cells.Item(66, 2).Value = 44.239699999999999
$ws.Cells.Item(67, 2).Value = 46.891500000000001
$ws.Cells.Item(68, 2).Value = 45.945
$ws.Cells.Item(69, 2).Value = 49.728999999999999
$ws.Cells.Item(70, 2).Value = 51.072800000000001
$ws.Cells.Item(71, 2).Value = 53.663699999999999
$ws.Cells.Item(72, 2).Value = 53.5533
$ws.Cells.Item(73, 2).Value = 53.7425
$ws.Cells.Item(74, 2).Value = 50.078800000000001
$ws.Cells.Item(75, 2).Value = 52.771299999999997
$ws.Cells.Item(76, 2).Value = 52.527299999999997
$ws.Cells.Item(77, 2).Value = 51.571100000000001
$ws.Cells.Item(78, 2).Value = 52.106200000000001
$ws.Cells.Item(79, 2).Value = 52.253700000000002
$ws.Cells.Item(80, 2).Value = 53.088099999999997
$ws.Cells.Item(81, 2).Value = 53.432200000000002
$ws.Cells.Item(82, 2).Value = 53.633499999999998
$ws.Cells.Item(83, 2).Value = 54.612000000000002
$ws.Cells.Item(84, 2).Value = 55.016599999999997
$ws.Cells.Item(85, 2).Value = 55.5413
$ws.Cells.Item(86, 2).Value = 53.408999999999999
$ws.Cells.Item(87, 2).Value = 55.884300000000003
$ws.Cells.Item(88, 2).Value = 55.821399999999997
$ws.Cells.Item(89, 2).Value = 56.468499999999999
$ws.Cells.Item(90, 2).Value = 56.987299999999998
$ws.Cells.Item(91, 2).Value = 56.502600000000001
$ws.Cells.Item(92, 2).Value = 56.5291
$ws.Cells.Item(93, 2).Value = 56.729799999999997
$ws.Cells.Item(94, 2).Value = 57.248699999999999
$ws.Cells.Item(95, 2).Value = 57.736400000000003
$ws.Cells.Item(96, 2).Value = 57.362900000000003
$ws.Cells.Item(97, 2).Value = 57.251199999999997
$ws.Cells.Item(98, 2).Value = 54.697800000000001
$ws.Cells.Item(99, 2).Value = 56.516599999999997
$ws.Cells.Item(100, 2).Value = 57.880099999999999
$ws.Cells.Item(101, 2).Value = 57.679099999999998
$ws.Cells.Item(102, 2).Value = 57.168700000000001
$ws.Cells.Item(103, 2).Value = 58.324316891000002
$ws.Cells.Item(104, 2).Value = 58.591552018068469
$ws.Cells.Item(105, 2).Value = 58.465340320980822
$ws.Cells.Item(106, 2).Value = 58.254994767996607
$ws.Cells.Item(107, 2).Value = 58.122121622072477
$ws.Cells.Item(108, 2).Value = 58.311300931422636
$ws.Cells.Item(109, 2).Value = 57.434632950310579
$ws.Cells.Item(110, 2).Value = 55.256828771374344
$ws.Cells.Item(111, 2).Value = 56.362688158707861
$ws.Cells.Item(112, 2).Value = 56.678325662885896
$ws.Cells.Item(113, 2).Value = 57.157353838496441
$ws.Cells.Item(114, 2).Value = 57.502189788770472
$ws.Cells.Item(115, 2).Value = 57.160908819571731
